# "Add Stage 3 (B1-B2) evidence"
#
# Fills in the TxHash evidence rows (A2/A3) on the "B1" and "B2" sheets,
# which previously held placeholder strings, moves the Info sheet's
# remembered selection off of D8 (and off of being the active tab), and
# makes "B2" the active/selected sheet with A4 selected on both B1 and B2.

$wb = $excel.ActiveWorkbook

# --- Info sheet: lose tabSelected, remember a different selection (B2) ---
$wsInfo = $wb.Worksheets.Item("Info")
[void]$wsInfo.Range("B2").Select()

# --- B1: fill in the two evidence hashes, remember selection A4 ---
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "2EAB50CBBC1B035073C1B33E9C34B048742633EA12287CCF9E8C75AECA263A9C"
$wsB1.Range("A3").Value = "75BD1A9C756EA7AA850073DB4786F953BF355D76549D508AA0B5CFBD1B6D8417"
[void]$wsB1.Range("A4").Select()

# --- B2: fill in the two evidence hashes, remember selection A4 ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "B05C79D763058DAF9BD11171E27FA968DC956355E97A05B98618B02BA77BB81E"
$wsB2.Range("A3").Value = "3A565BEB74FC8FA693CA3C245E355A01927E46A0A6D764BD6F26559C97704DAE"
[void]$wsB2.Range("A4").Select()

# B2 becomes the active sheet (tabSelected="1" / workbook activeTab points at it)
[void]$wsB2.Activate()
